$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registrationinfo")

$ws.Range("C2").Value = "a@mail.com"
$ws.Range("C3").Value = "b@mail.com"
$ws.Range("C4").Value = "c@mail.com"
$ws.Range("C5").Value = "d@mail.com"
$ws.Range("C6").Value = "e@mail.com"

$ws.Range("C7").Select()
